# "break out stock.yaml completed"
# The "day" sheet gets a fresh batch of scraped rows (366-380) appended,
# duplicating the previous batch (rows 351-365) with an updated timestamp
# and a handful of refreshed numeric values. The old batch's bsecode
# column (D351:D365) is normalized from text to a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1. Duplicate rows 351:365 into new rows 366:380 (new scrape batch) ---
#        (done before the D-column normalization below so the new batch
#        keeps its original text-typed bsecode values)
$src = $ws.Range("A351:I365")
$dst = $ws.Range("A366:I380")
$src.Copy($dst)

# --- 2. Normalize D351:D365 (bsecode) from text to numeric values ---
for ($r = 351; $r -le 365; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $numVal = [double]($cell.Value2)
    $cell.Value = $numVal
}

# --- 3. Refresh the handful of cells whose values changed in the new batch ---
$ws.Range("G367").Value = 56724
$ws.Range("G368").Value = 340278
$ws.Range("E372").Value = -1.7
$ws.Range("F372").Value = 1715.45
$ws.Range("G373").Value = 3453252
$ws.Range("G374").Value = 3121337
$ws.Range("E375").Value = -0.4
$ws.Range("F375").Value = 538.5
$ws.Range("G375").Value = 478064
$ws.Range("G378").Value = 4645867
$ws.Range("G380").Value = 15858298

# --- 4. Update the "Date Time" column for the new batch ---
$ws.Range("I366:I380").Value = "15/08/2024 11:35:01"
